$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 68.75
$ws.Range("I6").Value = 68.75
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 206.25
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -94.25
$ws.Range("N6").ClearContents()
$ws.Range("H99").Value = 2313.5
$ws.Range("J99").Value = 2993.3333
$ws.Range("L99").Value = 8979.999899999999
$ws.Range("N99").Value = -11975.9999
$ws.Range("H132").Value = 1180.375
$ws.Range("I132").Value = 1089.4193
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3268.2579
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -738.2579000000001
$ws.Range("N132").Value = -17060
$ws.Range("H135").Value = 647.38464
$ws.Range("I135").Value = 271.8889
$ws.Range("J135").Value = 1492.25
$ws.Range("K135").Value = 2447.0001
$ws.Range("L135").Value = 13430.25
$ws.Range("M135").Value = 87.99990000000025
$ws.Range("N135").Value = -18500.25
$ws.Range("H137").Value = 1787.5
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 1828.5714
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 5485.7142
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -10585.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 615.6667
$ws.Range("I2").Value = 619
$ws.Range("K2").Value = 619
$ws.Range("M2").Value = -506
$ws.Range("H32").Value = 3572.551
$ws.Range("I32").Value = 2242.861
$ws.Range("K32").Value = 2242.861
$ws.Range("M32").Value = -1955.861
$ws.Range("H45").Value = 3649
$ws.Range("I45").Value = 800
$ws.Range("J45").Value = 4598.6665
$ws.Range("K45").Value = 800
$ws.Range("L45").Value = 4598.6665
$ws.Range("M45").Value = -423
$ws.Range("N45").Value = -5352.6665
$ws.Range("H74").Value = 4460.0713
$ws.Range("I74").Value = 4461.6665
$ws.Range("K74").Value = 4461.6665
$ws.Range("M74").Value = -3587.6665
$ws.Range("H77").Value = 4460.0713
$ws.Range("I77").Value = 4461.6665
$ws.Range("K77").Value = 22308.3325
$ws.Range("M77").Value = -17940.3325
$ws.Range("H110").Value = 2896.6667
$ws.Range("I110").Value = 2476
$ws.Range("K110").Value = 2476
$ws.Range("M110").Value = -431
$ws.Range("H116").Value = 615.6667
$ws.Range("I116").Value = 619
$ws.Range("K116").Value = 619
$ws.Range("M116").Value = 1675
$ws.Range("H122").Value = 2357
$ws.Range("I122").Value = 2410
$ws.Range("K122").Value = 7230
$ws.Range("M122").Value = -4780
$ws.Range("H132").Value = 2074.862
$ws.Range("I132").Value = 1224.6471
$ws.Range("J132").Value = 3279.3333
$ws.Range("K132").Value = 3673.9413
$ws.Range("L132").Value = 9837.999899999999
$ws.Range("M132").Value = -1143.9413
$ws.Range("N132").Value = -14897.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 615.6667
$ws.Range("I3").Value = 619
$ws.Range("K3").Value = 619
$ws.Range("M3").Value = -505
$ws.Range("H107").Value = 2479.3333
$ws.Range("I107").Value = 1564.5454
$ws.Range("J107").Value = 4995
$ws.Range("K107").Value = 1564.5454
$ws.Range("L107").Value = 4995
$ws.Range("M107").Value = 355.4546
$ws.Range("N107").Value = -8835
$ws.Range("H134").Value = 8110.9653
$ws.Range("I134").Value = 9398.096
$ws.Range("J134").Value = 4732.25
$ws.Range("K134").Value = 28194.288
$ws.Range("L134").Value = 14196.75
$ws.Range("M134").Value = -25659.288
$ws.Range("N134").Value = -19266.75
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 870
$ws.Range("I16").Value = 826.6667
$ws.Range("K16").Value = 826.6667
$ws.Range("M16").Value = -539.6667
$ws.Range("H31").Value = 2173.9119
$ws.Range("I31").Value = 957.48
$ws.Range("J31").Value = 5552.8887
$ws.Range("K31").Value = 957.48
$ws.Range("L31").Value = 5552.8887
$ws.Range("M31").Value = -662.48
$ws.Range("N31").Value = -6142.8887
$ws.Range("H34").Value = 2173.9119
$ws.Range("I34").Value = 957.48
$ws.Range("J34").Value = 5552.8887
$ws.Range("K34").Value = 957.48
$ws.Range("L34").Value = 5552.8887
$ws.Range("M34").Value = -755.48
$ws.Range("N34").Value = -5956.8887
$ws.Range("H58").Value = 1396.2727
$ws.Range("I58").Value = 1307.5
$ws.Range("J58").Value = 1633
$ws.Range("K58").Value = 1307.5
$ws.Range("L58").Value = 1633
$ws.Range("M58").Value = -1104.5
$ws.Range("N58").Value = -2039
$ws.Range("H105").Value = 984.8
$ws.Range("I105").Value = 984.8
$ws.Range("K105").Value = 984.8
$ws.Range("M105").Value = 762.2
$ws.Range("H113").Value = 870
$ws.Range("I113").Value = 826.6667
$ws.Range("K113").Value = 826.6667
$ws.Range("M113").Value = 1343.3333
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -64910
$ws.Range("H132").Value = 2009.9429
$ws.Range("I132").Value = 924.2381
$ws.Range("J132").Value = 3638.5
$ws.Range("K132").Value = 2772.7143
$ws.Range("L132").Value = 10915.5
$ws.Range("M132").Value = -242.7143000000001
$ws.Range("N132").Value = -15975.5
$ws.Range("H134").Value = 1030.091
$ws.Range("I134").Value = 1013.1
$ws.Range("K134").Value = 3039.3
$ws.Range("M134").Value = -504.3000000000002
$ws.Range("H136").Value = 1396.2727
$ws.Range("I136").Value = 1307.5
$ws.Range("J136").Value = 1633
$ws.Range("K136").Value = 3922.5
$ws.Range("L136").Value = 4899
$ws.Range("M136").Value = -1372.5
$ws.Range("N136").Value = -9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 443
$ws.Range("I7").Value = 239
$ws.Range("J7").Value = 715
$ws.Range("K7").Value = 717
$ws.Range("L7").Value = 2145
$ws.Range("M7").Value = -605
$ws.Range("N7").Value = -2369
$ws.Range("H115").Value = 4979.4
$ws.Range("J115").Value = 5724.25
$ws.Range("L115").Value = 17172.75
$ws.Range("N115").Value = -19522.75
$ws.Range("H122").Value = 1082.0834
$ws.Range("I122").Value = 883.3333
$ws.Range("J122").Value = 1280.8334
$ws.Range("K122").Value = 7949.9997
$ws.Range("L122").Value = 11527.5006
$ws.Range("M122").Value = -5499.9997
$ws.Range("N122").Value = -16427.5006
$ws.Range("H131").Value = 11922108
$ws.Range("J131").Value = 20706.629
$ws.Range("L131").Value = 62119.887
$ws.Range("N131").Value = -72199.887
$ws.Range("H132").Value = 1944.4445
$ws.Range("I132").Value = 1585.7142
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 14271.4278
$ws.Range("L132").Value = 28800
$ws.Range("M132").Value = -11741.4278
$ws.Range("N132").Value = -33860
$ws.Range("H138").Value = 2036.1818
$ws.Range("I138").Value = 1628.2858
$ws.Range("K138").Value = 4884.857400000001
$ws.Range("M138").Value = 255.1425999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2130.4285
$ws.Range("I102").Value = 2347.3333
$ws.Range("K102").Value = 2347.3333
$ws.Range("M102").Value = -725.3332999999998
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H132").Value = 3853.875
$ws.Range("I132").Value = 2569.2856
$ws.Range("J132").Value = 5652.3
$ws.Range("K132").Value = 7707.8568
$ws.Range("L132").Value = 16956.9
$ws.Range("M132").Value = -5177.8568
$ws.Range("N132").Value = -22016.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1499
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29703.334
$ws.Range("J70").Value = 29703.334
$ws.Range("L70").Value = 29703.334
$ws.Range("N70").Value = -30333.334
$ws.Range("H73").Value = 29703.334
$ws.Range("J73").Value = 29703.334
$ws.Range("L73").Value = 29703.334
$ws.Range("N73").Value = -31887.334
$ws.Range("H136").Value = 3661.5
$ws.Range("I136").Value = 3435.6428
$ws.Range("J136").Value = 3925
$ws.Range("K136").Value = 10306.9284
$ws.Range("L136").Value = 11775
$ws.Range("M136").Value = -7756.928400000001
